# Add "Model Metrics on Validation Data" table (rows 51-57) to the Data sheet.
# Values / write order chosen to match the order new strings were
# originally interned into the workbook's shared-string table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- column headers for the metrics table (written first) ---
$ws.Range("C52").Value = " recall"
$ws.Range("D52").Value = " precision"
$ws.Range("E52").Value = " false negatives"

# --- row labels ---
$ws.Range("A53").Value = "Logistic"
$ws.Range("A54").Value = "1 Hidden Layer"
$ws.Range("A55").Value = "3 Hidden Layer"
$ws.Range("A56").Value = "5 Hidden Layer"
$ws.Range("A57").Value = "10 Hidden Layer"

# --- section title ---
$ws.Range("A51").Value = "Model Metrics on Validation Data"

# --- remaining headers ---
$ws.Range("B52").Value = "accuracy"
$ws.Range("G52").Value = "False Positive Rate"

# --- data rows ---
$ws.Range("B53").Value = 6.6613262272682503
$ws.Range("C53").Value = 0.58123249299719804
$ws.Range("D53").Value = 0.238841794168247
$ws.Range("E53").Value = 0.781752488158998
$ws.Range("F53").Value = 2736.99775910364
$ws.Range("G53").Formula = "=1-C53"

$ws.Range("B54").Value = 0.69020196750384399
$ws.Range("C54").Value = 0.69103641456582598
$ws.Range("D54").Value = 0.490495910697958
$ws.Range("E54").Value = 0.83313285730132203
$ws.Range("F54").Value = 1828.5030812324901
$ws.Range("G54:G57").Formula = "=1-C54"

$ws.Range("B55").Value = 0.53603627344473403
$ws.Range("C55").Value = 0.73732492997198795
$ws.Range("D55").Value = 0.703700856234179
$ws.Range("E55").Value = 0.760432809841733
$ws.Range("F55").Value = 1065.7400560224

$ws.Range("B56").Value = 0.53103721742870402
$ws.Range("C56").Value = 0.73452380952380902
$ws.Range("D56").Value = 0.72523372303537903
$ws.Range("E56").Value = 0.742979875098423
$ws.Range("F56").Value = 985.60952380952301

$ws.Range("B57").Value = 0.45037841686681501
$ws.Range("C57").Value = 0.79782913165266101
$ws.Range("D57").Value = 0.77460099241646696
$ws.Range("E57").Value = 0.81370940976450001
$ws.Range("F57").Value = 810.21960784313706

# --- restore the view state reported in the target workbook ---
$ws.Range("G53:G57").Select()
